$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 34
$ws.Range("B3").Value = 417640
$ws.Range("B4").Value = 4760
$ws.Range("B5").Value = 4760
